# ToDo list update: re-prioritise tasks, drop two finished items, add four
# new ones, and re-point the selection at B5 (the "Vertex welding" estimate).
#
# Shared-string order in the saved workbook follows the order in which new
# strings are first written, so the two "Refactoring" / "Error handling"
# rows are written before the "Make it spin" / "Data drive lighting" /
# "Fix sorting issues" rows to match the target sharedStrings.xml layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 6-7: new tasks (written first so their strings land earlier in sst)
$ws.Range("A6").Value = "Refactoring - we need consistency across the board"
$ws.Range("B6").Value = 21
$ws.Range("A7").Value = "Error handling - go on a robustness run, set standards for future work"
$ws.Range("B7").Value = 21

# Rows 2-4: new tasks inserted ahead of the carried-over rows
$ws.Range("A2").Value = "Make it spin"
$ws.Range("B2").Value = 3
$ws.Range("A3").Value = "Data drive lighting"
$ws.Range("B3").Value = 3
$ws.Range("A4").Value = "Fix sorting issues"
$ws.Range("B4").Value = 14

# Row 5: carried-over task, now with an updated estimate
$ws.Range("A5").Value = "Vertex welding in model compiler"
$ws.Range("B5").Value = 14

# Row 8: carried-over task, moved to the bottom with an updated estimate
$ws.Range("A8").Value = "Textured surfaces"
$ws.Range("B8").Value = 35

# Match the workbook's saved selection state
$ws.Range("B5").Select()
